$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(4, 7).Value = 3.6
$ws.Cells.Item(4, 8).Value = 2.52
$ws.Cells.Item(4, 9).Value = 2.45
$ws.Cells.Item(4, 10).Value = 4.35
$ws.Cells.Item(4, 11).Value = 1.78
$ws.Cells.Item(4, 12).Value = 3.15
$ws.Cells.Item(4, 13).Value = 1.16
$ws.Cells.Item(4, 14).Value = 4.55
$ws.Cells.Item(4, 17).Value = 2.82
$ws.Cells.Item(4, 25).Value = 2.15
$ws.Cells.Item(4, 26).Value = 1.62
$ws.Cells.Item(4, 27).Value = 7.2
$ws.Cells.Item(4, 28).Value = 17.5
$ws.Cells.Item(4, 30).Value = 60
$ws.Cells.Item(4, 31).Value = 45
$ws.Cells.Item(4, 32).Value = 65
$ws.Cells.Item(4, 33).Value = 4.55
$ws.Cells.Item(4, 34).Value = 5.2
$ws.Cells.Item(4, 35).Value = 17.5
$ws.Cells.Item(4, 38).Value = 5.6
$ws.Cells.Item(4, 39).Value = 10.5
$ws.Cells.Item(4, 40).Value = 10
$ws.Cells.Item(4, 41).Value = 28
$ws.Cells.Item(4, 42).Value = 26
$ws.Cells.Item(4, 43).Value = 45
$ws.Cells.Item(5, 7).Value = 3.75
$ws.Cells.Item(5, 8).Value = 2.72
$ws.Cells.Item(5, 9).Value = 2.25
$ws.Cells.Item(5, 10).Value = 4.4
$ws.Cells.Item(5, 11).Value = 1.85
$ws.Cells.Item(5, 12).Value = 2.92
$ws.Cells.Item(5, 13).Value = 1.14
$ws.Cells.Item(5, 14).Value = 5
$ws.Cells.Item(5, 16).Value = 2.3
$ws.Cells.Item(5, 17).Value = 2.62
$ws.Cells.Item(5, 21).Value = 4.75
$ws.Cells.Item(5, 25).Value = 2.1
$ws.Cells.Item(5, 26).Value = 1.65
$ws.Cells.Item(5, 27).Value = 7.8
$ws.Cells.Item(5, 28).Value = 18.5
$ws.Cells.Item(5, 29).Value = 13.5
$ws.Cells.Item(5, 30).Value = 60
$ws.Cells.Item(5, 31).Value = 45
$ws.Cells.Item(5, 32).Value = 60
$ws.Cells.Item(5, 33).Value = 5
$ws.Cells.Item(5, 34).Value = 5.5
$ws.Cells.Item(5, 35).Value = 17.5
$ws.Cells.Item(5, 38).Value = 5.6
$ws.Cells.Item(5, 39).Value = 9.5
$ws.Cells.Item(5, 40).Value = 9.5
$ws.Cells.Item(5, 41).Value = 23
$ws.Cells.Item(5, 42).Value = 23
$ws.Cells.Item(7, 17).Value = 2.2
$ws.Cells.Item(7, 18).Value = 1.67
$ws.Cells.Item(8, 7).Value = 1.3
$ws.Cells.Item(8, 8).Value = 5
$ws.Cells.Item(8, 9).Value = 9.5
$ws.Cells.Item(8, 10).Value = 1.83
$ws.Cells.Item(8, 11).Value = 2.3
$ws.Cells.Item(8, 12).Value = 9.5
$ws.Cells.Item(8, 17).Value = 2
$ws.Cells.Item(8, 18).Value = 1.85
$ws.Cells.Item(8, 21).Value = 3.5
$ws.Cells.Item(8, 22).Value = 1.3
$ws.Cells.Item(8, 25).Value = 2.5
$ws.Cells.Item(8, 26).Value = 1.5
$ws.Cells.Item(8, 29).Value = 9
$ws.Cells.Item(9, 7).Value = 1.7
$ws.Cells.Item(9, 33).Value = 11
$ws.Cells.Item(9, 40).Value = 15
$ws.Cells.Item(10, 7).Value = 2.25
$ws.Cells.Item(10, 9).Value = 3.2
$ws.Cells.Item(10, 10).Value = 3
$ws.Cells.Item(10, 14).Value = 8.5
$ws.Cells.Item(10, 27).Value = 7
$ws.Cells.Item(10, 28).Value = 10
$ws.Cells.Item(10, 31).Value = 19
$ws.Cells.Item(10, 40).Value = 12
$ws.Cells.Item(10, 43).Value = 41
$ws.Cells.Item(11, 17).Value = 2.1
$ws.Cells.Item(11, 18).Value = 1.73
$ws.Cells.Item(11, 25).Value = 1.83
$ws.Cells.Item(11, 26).Value = 1.83
$ws.Cells.Item(12, 13).Value = 1.05
$ws.Cells.Item(12, 14).Value = 11
$ws.Cells.Item(12, 17).Value = 1.88
$ws.Cells.Item(12, 18).Value = 1.98
$ws.Cells.Item(12, 25).Value = 1.73
$ws.Cells.Item(13, 15).Value = 1.25
$ws.Cells.Item(13, 16).Value = 3.75
$ws.Cells.Item(13, 17).Value = 1.8
$ws.Cells.Item(13, 18).Value = 2
$ws.Cells.Item(13, 21).Value = 3
$ws.Cells.Item(13, 22).Value = 1.36
$ws.Cells.Item(13, 25).Value = 1.67
$ws.Cells.Item(14, 7).Value = 1.42
$ws.Cells.Item(14, 10).Value = 1.91
$ws.Cells.Item(14, 12).Value = 6
$ws.Cells.Item(14, 17).Value = 1.48
$ws.Cells.Item(14, 18).Value = 2.6
$ws.Cells.Item(14, 19).Value = 1.8
$ws.Cells.Item(14, 20).Value = 2.05
$ws.Cells.Item(14, 25).Value = 1.67
$ws.Cells.Item(14, 26).Value = 2.1
$ws.Cells.Item(14, 32).Value = 21
$ws.Cells.Item(15, 14).Value = 13
$ws.Cells.Item(15, 17).Value = 1.75
$ws.Cells.Item(15, 25).Value = 1.57
$ws.Cells.Item(16, 15).Value = 1.29
$ws.Cells.Item(16, 16).Value = 3.5
$ws.Cells.Item(16, 17).Value = 1.9
$ws.Cells.Item(16, 18).Value = 1.95
$ws.Cells.Item(16, 25).Value = 1.67
$ws.Cells.Item(17, 9).Value = 1.73
$ws.Cells.Item(17, 10).Value = 4.5
$ws.Cells.Item(17, 25).Value = 1.62
$ws.Cells.Item(17, 29).Value = 15
$ws.Cells.Item(17, 32).Value = 34
$ws.Cells.Item(17, 39).Value = 9.5
$ws.Cells.Item(20, 7).Value = 1.42
$ws.Cells.Item(20, 8).Value = 4
$ws.Cells.Item(20, 9).Value = 6.5
$ws.Cells.Item(20, 10).Value = 1.88
$ws.Cells.Item(20, 11).Value = 2.32
$ws.Cells.Item(20, 12).Value = 6
$ws.Cells.Item(20, 15).Value = 1.18
$ws.Cells.Item(20, 16).Value = 4.3
$ws.Cells.Item(20, 17).Value = 1.57
$ws.Cells.Item(20, 18).Value = 2.1
$ws.Cells.Item(20, 21).Value = 2.45
$ws.Cells.Item(20, 22).Value = 1.49
$ws.Cells.Item(20, 25).Value = 1.8
$ws.Cells.Item(20, 26).Value = 1.91
$ws.Cells.Item(20, 27).Value = 6.6
$ws.Cells.Item(20, 28).Value = 6.3
$ws.Cells.Item(20, 29).Value = 6.8
$ws.Cells.Item(20, 30).Value = 8.5
$ws.Cells.Item(20, 31).Value = 9
$ws.Cells.Item(20, 32).Value = 17.5
$ws.Cells.Item(20, 33).Value = 12.5
$ws.Cells.Item(20, 34).Value = 7.1
$ws.Cells.Item(20, 35).Value = 13
$ws.Cells.Item(20, 36).Value = 50
$ws.Cells.Item(20, 37).Value = 350
$ws.Cells.Item(20, 38).Value = 16
$ws.Cells.Item(20, 39).Value = 35
$ws.Cells.Item(20, 40).Value = 16.5
$ws.Cells.Item(20, 41).Value = 110
$ws.Cells.Item(20, 42).Value = 50
$ws.Cells.Item(20, 43).Value = 40
$ws.Cells.Item(21, 7).Value = 2.22
$ws.Cells.Item(21, 8).Value = 3.15
$ws.Cells.Item(21, 9).Value = 2.92
$ws.Cells.Item(21, 10).Value = 2.85
$ws.Cells.Item(21, 11).Value = 2.02
$ws.Cells.Item(21, 12).Value = 3.5
$ws.Cells.Item(21, 15).Value = 1.27
$ws.Cells.Item(21, 16).Value = 3.4
$ws.Cells.Item(21, 17).Value = 1.83
$ws.Cells.Item(21, 18).Value = 1.78
$ws.Cells.Item(21, 21).Value = 3.1
$ws.Cells.Item(21, 22).Value = 1.32
$ws.Cells.Item(21, 23).Value = 1.38
$ws.Cells.Item(21, 24).Value = 2.45
$ws.Cells.Item(21, 25).Value = 1.71
$ws.Cells.Item(21, 26).Value = 2.03
$ws.Cells.Item(21, 27).Value = 6.7
$ws.Cells.Item(21, 28).Value = 9.25
$ws.Cells.Item(21, 29).Value = 7.5
$ws.Cells.Item(21, 30).Value = 18
$ws.Cells.Item(21, 31).Value = 14.5
$ws.Cells.Item(21, 32).Value = 21
$ws.Cells.Item(21, 33).Value = 9.5
$ws.Cells.Item(21, 34).Value = 5.5
$ws.Cells.Item(21, 35).Value = 10.75
$ws.Cells.Item(21, 36).Value = 40
$ws.Cells.Item(21, 37).Value = 250
$ws.Cells.Item(21, 38).Value = 8.25
$ws.Cells.Item(21, 39).Value = 13
$ws.Cells.Item(21, 40).Value = 8.75
$ws.Cells.Item(21, 41).Value = 29
$ws.Cells.Item(21, 42).Value = 19.5
$ws.Cells.Item(21, 43).Value = 24
$ws.Cells.Item(23, 7).Value = 2.25
$ws.Cells.Item(23, 8).Value = 3.2
$ws.Cells.Item(23, 9).Value = 3.1
$ws.Cells.Item(23, 10).Value = 2.95
$ws.Cells.Item(23, 11).Value = 2
$ws.Cells.Item(23, 12).Value = 3.7
$ws.Cells.Item(23, 13).Value = 1.08
$ws.Cells.Item(23, 14).Value = 6.6
$ws.Cells.Item(23, 15).Value = 1.37
$ws.Cells.Item(23, 16).Value = 2.85
$ws.Cells.Item(23, 17).Value = 2.1
$ws.Cells.Item(23, 18).Value = 1.65
$ws.Cells.Item(23, 21).Value = 3.55
$ws.Cells.Item(23, 22).Value = 1.25
$ws.Cells.Item(23, 23).Value = 1.47
$ws.Cells.Item(23, 24).Value = 2.5
$ws.Cells.Item(23, 25).Value = 1.85
$ws.Cells.Item(23, 26).Value = 1.85
$ws.Cells.Item(23, 27).Value = 6.9
$ws.Cells.Item(23, 28).Value = 10.25
$ws.Cells.Item(23, 29).Value = 9.25
$ws.Cells.Item(23, 30).Value = 22
$ws.Cells.Item(23, 31).Value = 19.5
$ws.Cells.Item(23, 32).Value = 32
$ws.Cells.Item(23, 33).Value = 6.6
$ws.Cells.Item(23, 34).Value = 6.1
$ws.Cells.Item(23, 35).Value = 15
$ws.Cells.Item(23, 36).Value = 75
$ws.Cells.Item(23, 38).Value = 8.75
$ws.Cells.Item(23, 39).Value = 15.5
$ws.Cells.Item(23, 40).Value = 11
$ws.Cells.Item(23, 41).Value = 40
$ws.Cells.Item(23, 42).Value = 28
$ws.Cells.Item(23, 43).Value = 37
